$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.003183900000000101
$ws.Range("E2").Value = 0.3723852725881063
$ws.Range("G2").Value = 0.2494892361374987
$ws.Range("I2").Value = 0.3669021
$ws.Range("L2").Value = 0.5961429402307629
$ws.Range("M2").Value = 0.08239116666666667
$ws.Range("N2").Value = 12.89781179131339
$ws.Range("O2").Value = 3.475337169561457

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.07396955163932978
$ws.Range("E2").Value = 0.3637656372628776
$ws.Range("I2").Value = 0.5780978354701471
$ws.Range("L2").Value = 0.2763239402451296
$ws.Range("M2").Value = 0.08143291666666665
$ws.Range("N2").Value = 9.06843726173757
$ws.Range("O2").Value = 3.840635583889667

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.05173369849322602
$ws.Range("B2").Value = 0.02941088209256308
$ws.Range("E2").Value = 0.1975555792241696
$ws.Range("I2").Value = 0.4424245787642941
$ws.Range("M2").Value = 0.0525461666666667
$ws.Range("N2").Value = 8.760925509998557
$ws.Range("O2").Value = 5.233712211049507

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 0.001164164422075072

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 0.1863654784114568
